$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '79.504.92'
$ws.Range('E2').Value = '  +4.14%  '
$ws.Range('D3').Value = '3.163.58'
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '206.58'
$ws.Range('E5').Value = '  +4.17%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '627.17'
$ws.Range('E6').Value = '  +1.96%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.268'
$ws.Range('E7').Value = '  +27.46%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.591'
$ws.Range('E9').Value = '  +7.24%  '
$ws.Range('D10').Value = '3.160.29'
$ws.Range('E10').Value = '  +2.48%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.596'
$ws.Range('E11').Value = '  +35.07%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000251'
$ws.Range('E12').Value = '  +29.45%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.164'
$ws.Range('E13').Value = '  +2.15%  '
$ws.Range('E14').Value = '  +0.94%  '
$ws.Range('D15').Value = '3.735.50'
$ws.Range('E15').Value = '  +2.43%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '31.50'
$ws.Range('E16').Value = '  +7.74%  '
$ws.Range('D17').Value = '79.608.50'
$ws.Range('E17').Value = '  +4.43%  '
$ws.Range('D18').Value = '3.159.64'
$ws.Range('E18').Value = '  +2.47%  '
$ws.Range('D19').Value = '14.33'
$ws.Range('E19').Value = '  +5.41%  '
$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '2.95'
$ws.Range('E20').Value = '  +14.72%  '
$ws.Range('B21').Value = 'BitcoinCash'
$ws.Range('C21').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D21').Value = '437.56'
$ws.Range('E21').Value = '  +14.77%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.12'
$ws.Range('E22').Value = '  -0.32%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.24'
$ws.Range('E23').Value = '  +18.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.77'
$ws.Range('E24').Value = '  +4.54%  '
$ws.Range('D25').Value = '3.327.43'
$ws.Range('E25').Value = '  +2.81%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '76.23'
$ws.Range('E26').Value = '  +5.45%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '4.68'
$ws.Range('E27').Value = '  +6.68%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.89'
$ws.Range('E28').Value = '  +9.00%  '
$ws.Range('D29').Value = '0.999'
$ws.Range('E29').Value = '  -0.32%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0000121'
$ws.Range('E30').Value = '  +11.20%  '
$ws.Range('B31').Value = 'Binance-PegBSC-USD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  +0.05%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '9.01'
$ws.Range('E32').Value = '  +8.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '551.02'
$ws.Range('E33').Value = '  +10.34%  '
$ws.Range('D34').Value = '1.48'
$ws.Range('E34').Value = '  +4.10%  '
$ws.Range('B35').Value = 'PancakeSwap'
$ws.Range('C35').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.00'
$ws.Range('E35').Value = '  +4.39%  '
$ws.Range('B36').Value = 'Kaspa'
$ws.Range('C36').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.148'
$ws.Range('E36').Value = '  +20.37%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '23.14'
$ws.Range('E37').Value = '  +11.27%  '
$ws.Range('E38').Value = '  +16.91%  '
$ws.Range('D39').Value = '0.998'
$ws.Range('E39').Value = '  -0.08%  '
$ws.Range('E40').Value = '  +7.39%  '
$ws.Range('D41').Value = '20.77'
$ws.Range('E41').Value = '  +3.52%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '163.24'
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.63'
$ws.Range('E43').Value = '  +9.92%  '
$ws.Range('E44').Value = '  -0.02%  '
$ws.Range('D45').Value = '188.09'
$ws.Range('E45').Value = '  -3.84%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.80'
$ws.Range('E46').Value = '  +9.32%  '
$ws.Range('E47').Value = '  +9.73%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.780'
$ws.Range('E48').Value = '  -2.62%  '
$ws.Range('B49').Value = 'ImmutableX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.30'
$ws.Range('E49').Value = '  +4.04%  '
$ws.Range('B50').Value = 'OKB'
$ws.Range('C50').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D50').Value = '43.32'
$ws.Range('E50').Value = '  +4.79%  '
$ws.Range('E51').Value = '  +9.19%  '
